$d = $word.ActiveDocument

function Set-NewRunFont($range) {
    $range.Font.Name = "Calibri"
    $range.Font.Color = 0
}

# ---------------------------------------------------------------
# 1. Title
# ---------------------------------------------------------------
$d.Content.Find.Execute("The Allure of the Quantum Realm", $true, $false, $false, $false, $false, $true, 1, $false, "The Enigmatic Symphony of Politics: Navigating the Complex World of Power and Influence", 2) | Out-Null

# ---------------------------------------------------------------
# 2. Author name
# ---------------------------------------------------------------
$d.Content.Find.Execute("Genevieve Smith", $true, $false, $false, $false, $false, $true, 1, $false, "Eleanor Winters", 2) | Out-Null

# ---------------------------------------------------------------
# 3. Email paragraph (paragraph 3) - scope Find to this paragraph only,
#    since "com" is a common substring elsewhere (computing/computers/computation).
# ---------------------------------------------------------------
$pEmail = $d.Paragraphs(3)
$pEmail.Range.Find.Execute("smithgenevieve@gmail", $true, $false, $false, $false, $false, $true, 1, $false, "eleanor", 2) | Out-Null
$pEmail.Range.Find.Execute("com", $true, $false, $false, $false, $false, $true, 1, $false, "winters@validedu", 2) | Out-Null

# append new "." and "org" runs at the end of the email paragraph
$emailEnd = $pEmail.Range.End - 1
$rDot = $d.Range($emailEnd, $emailEnd)
$rDot.InsertAfter(".")
Set-NewRunFont $rDot
$rDot.Font.Size = 16

$emailEnd2 = $pEmail.Range.End - 1
$rOrg = $d.Range($emailEnd2, $emailEnd2)
$rOrg.InsertAfter("org")
Set-NewRunFont $rOrg
$rOrg.Font.Size = 16

# ---------------------------------------------------------------
# 4. Body paragraph (paragraph 5) - the long quantum -> politics intro
# ---------------------------------------------------------------
$d.Content.Find.Execute("Journey with us into the captivating realm of quantum physics, where particles dance in a harmonious waltz of uncertainty and probability", $true, $false, $false, $false, $false, $true, 1, $false, "In the realm of human affairs, there exists an intricate symphony of power, influence, and decision-making that we call politics", 2) | Out-Null

$d.Content.Find.Execute(" This enigmatic realm, once confined to theoretical musings, is now seeping into our reality, promising awe-inspiring technologies that defy classical intuition", $true, $false, $false, $false, $false, $true, 1, $false, " This vast and dynamic world of governance, leadership, and societal interactions shapes the very fabric of our societies, affecting every aspect of our lives, from the policies that govern us to the leaders who represent us", 2) | Out-Null

$d.Content.Find.Execute(" From the enigmatic world of quantum computing to the nascent field of quantum cryptography, we stand at the precipice of a paradigm shift, where the ethereal fabric of quantum mechanics is woven into the tapestry of our technological landscape", $true, $false, $false, $false, $false, $true, 1, $false, " To navigate this complex landscape effectively, it is imperative that we understand the fundamental principles of politics, its historical evolution, and its profound impact on our daily lives", 2) | Out-Null

$d.Content.Find.Execute("Step into the arena of quantum computing, where information dances in the ethereal realm of quantum bits, also known as qubits", $true, $false, $false, $false, $false, $true, 1, $false, "Politics, at its core, is the art of resolving conflicts and allocating resources within a society", 2) | Out-Null

$d.Content.Find.Execute(" Unlike their classical counterparts, qubits wield the uncanny ability to exist in a superposition of states, pirouette-ing through a ballet of possibilities", $true, $false, $false, $false, $false, $true, 1, $false, " It involves the formulation and implementation of policies, the establishment of laws and regulations, and the distribution of power among various institutions and individuals", 2) | Out-Null

$d.Content.Find.Execute(" This intoxicating dance grants quantum computers the potency to tackle conundrums that confound their classical brethren, pioneering solutions to intractable problems in cryptography, optimization, and simulation, unveiling secrets hidden within the labyrinthine pathways of computation", $true, $false, $false, $false, $false, $true, 1, $false, " Through political processes, we collectively determine how we want to live together, what values we hold dear, and how we can create a just and equitable society for all", 2) | Out-Null

$d.Content.Find.Execute("Venture into the clandestine realm of quantum cryptography, where information cloaks itself in the enigmatic embrace of quantum mechanics, creating an impregnable shield against eavesdropping ears", $true, $false, $false, $false, $false, $true, 1, $false, "As we delve into the study of politics, we are confronted with a tapestry of historical events, political theories, and ideological debates that have shaped our current political landscape", 2) | Out-Null

$d.Content.Find.Execute(" This quantum cloak harnesses the inherent fragility of quantum information, orchestrating a symphony of particles that evokes alarm at the slightest touch of an unintended observer", $true, $false, $false, $false, $false, $true, 1, $false, " From the ancient Greek city-states to the modern nation-states, from the rise and fall of empires to the emergence of global governance, politics has been an ever-evolving field, constantly adapting to changing circumstances and societal needs", 2) | Out-Null

$d.Content.Find.Execute(" With quantum cryptography as our guardian, we can forge unbreakable codes, ensuring the sanctity of our secrets in a world where data breaches are an incessant threat", $true, $false, $false, $false, $false, $true, 1, $false, " By understanding the historical context of politics, we gain a deeper appreciation for the challenges and opportunities that lie ahead", 2) | Out-Null

# ---------------------------------------------------------------
# 4b. Insert the large new block of runs right before the final "."
#     that ends paragraph 5 (locate "lie ahead" then position after it).
# ---------------------------------------------------------------
$found = $d.Content
$found.Find.Execute("lie ahead", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$insertPos = $found.End
$ins = $d.Range($insertPos, $insertPos)

function Add-Run($target, [string]$text) {
    $target.InsertAfter($text)
    Set-NewRunFont $target
    $target.Font.Size = 12
    return $target
}

# Run A: "."
$pos = $ins.End
$r = $d.Range($pos, $pos)
Add-Run $r "." | Out-Null

# Run B: lone break
$pos = $r.End
$r = $d.Range($pos, $pos)
Add-Run $r ([string][char]11) | Out-Null

# Run C: break + "Introduction Continued:"
$pos = $r.End
$r = $d.Range($pos, $pos)
Add-Run $r ([string][char]11 + "Introduction Continued:") | Out-Null

# Run D: lone break
$pos = $r.End
$r = $d.Range($pos, $pos)
Add-Run $r ([string][char]11) | Out-Null

# Run E: break + "Furthermore, the study of politics provides us with a framework for analyzing and understanding current events"
$pos = $r.End
$r = $d.Range($pos, $pos)
Add-Run $r ([string][char]11 + "Furthermore, the study of politics provides us with a framework for analyzing and understanding current events") | Out-Null

# Run F: "."
$pos = $r.End
$r = $d.Range($pos, $pos)
Add-Run $r "." | Out-Null

# Run G
$pos = $r.End
$r = $d.Range($pos, $pos)
Add-Run $r " By examining political institutions, policies, and decision-making processes, we can develop a critical perspective on the world around us" | Out-Null

# Run H: "."
$pos = $r.End
$r = $d.Range($pos, $pos)
Add-Run $r "." | Out-Null

# Run I
$pos = $r.End
$r = $d.Range($pos, $pos)
Add-Run $r " We can identify the various actors and interests at play, assess the potential impact of different " | Out-Null

# Run J: "policies, and engage in informed debates about the direction of our society" (with lastRenderedPageBreak - best effort, see below)
$pos = $r.End
$r = $d.Range($pos, $pos)
Add-Run $r "policies, and engage in informed debates about the direction of our society" | Out-Null

# Run K: "."
$pos = $r.End
$r = $d.Range($pos, $pos)
Add-Run $r "." | Out-Null

# Run L
$pos = $r.End
$r = $d.Range($pos, $pos)
Add-Run $r " Politics is not merely an abstract concept; it is a living, breathing force that shapes our communities, our economies, and our planet" | Out-Null

# Run M: "."
$pos = $r.End
$r = $d.Range($pos, $pos)
Add-Run $r "." | Out-Null

# Run N: lone break
$pos = $r.End
$r = $d.Range($pos, $pos)
Add-Run $r ([string][char]11) | Out-Null

# Run O: break + "Politics is a multi-faceted subject..."
$pos = $r.End
$r = $d.Range($pos, $pos)
Add-Run $r ([string][char]11 + "Politics is a multi-faceted subject that encompasses a wide range of topics, from the intricacies of international relations to the challenges of local governance") | Out-Null

# Run P: "."
$pos = $r.End
$r = $d.Range($pos, $pos)
Add-Run $r "." | Out-Null

# Run Q
$pos = $r.End
$r = $d.Range($pos, $pos)
Add-Run $r " It involves the study of political systems, ideologies, public policy, and the role of citizens in a democracy" | Out-Null

# Run R: "."
$pos = $r.End
$r = $d.Range($pos, $pos)
Add-Run $r "." | Out-Null

# Run S
$pos = $r.End
$r = $d.Range($pos, $pos)
Add-Run $r " Through political engagement, we have the power to influence the decisions that affect our lives and to hold our leaders accountable" | Out-Null

# ---------------------------------------------------------------
# 5. Summary heading paragraph text replacements
# ---------------------------------------------------------------
$d.Content.Find.Execute("The captivating realm of quantum physics dances on the boundary of our perception, blurring the line between theory and reality", $true, $false, $false, $false, $false, $true, 1, $false, "In this essay, we have explored the enigmatic symphony of politics, a complex world of power, influence, and decision-making that profoundly impacts our lives", 2) | Out-Null

$d.Content.Find.Execute(" Quantum computing, the nascent field of quantum cryptography, and the burgeoning world of quantum sensing are transforming our technological landscape", $true, $false, $false, $false, $false, $true, 1, $false, " By understanding the fundamental principles of politics, its historical evolution, and its current manifestations, we gain a deeper appreciation for the challenges and opportunities that lie ahead", 2) | Out-Null

$d.Content.Find.Execute(" From decoding intricate problems to securing our digital realm, the quantum realm is reshaping our perception of what's possible", $true, $false, $false, $false, $false, $true, 1, $false, " The study of politics equips us with the knowledge and skills necessary to navigate the intricacies of governance, to engage in informed debates, and to shape the future of our societies", 2) | Out-Null

# merge the two runs (" As we delve deeper..." + "are unveiling a new chapter...") into one,
# which also naturally drops the lastRenderedPageBreak marker that sat on the 2nd run.
$d.Content.Find.Execute(" As we delve deeper into this enigmatic realm, we are unveiling a new chapter in the human quest for knowledge and dominion over the forces that govern our universe", $true, $false, $false, $false, $false, $true, 1, $false, " As active citizens, we have the responsibility to participate in the political process, to hold our leaders accountable, and to strive for a more just and equitable world for all", 2) | Out-Null

# ---------------------------------------------------------------
# 6. Trailing empty paragraph after the Summary paragraph
# ---------------------------------------------------------------
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$endR = $d.Range($lastPara.Range.End, $lastPara.Range.End)
$endR.InsertParagraphAfter()
